# Weekly update: insert 3 new rows of Espárragos price data (row 97-99)
# at the top of the Mercado Mayorista Lo Valledor de Santiago series,
# pushing all existing rows (old 97-179) down by 3 (new 100-182).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 97..99, shifting rows 97:179 down to 100:182.
$ws.Range("A97:R99").Insert()

# --- New row 97: Banquete ---
$ws.Cells.Item(97,1).Value = 6
$ws.Cells.Item(97,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(97,3).Value = "Metropolitana"
$ws.Cells.Item(97,4).Value = 45225
$ws.Cells.Item(97,5).Value = 13
$ws.Cells.Item(97,6).Value = 300000000
$ws.Cells.Item(97,7).Value = "Espárragos"
$ws.Cells.Item(97,8).Value = "Sin especificar"
$ws.Cells.Item(97,9).Value = "Banquete"
$ws.Cells.Item(97,10).Value = 4000
$ws.Cells.Item(97,11).Value = 1500
$ws.Cells.Item(97,12).Value = 1600
$ws.Cells.Item(97,13).Value = 1550
$ws.Cells.Item(97,14).Value = "$/kilo"
$ws.Cells.Item(97,15).Value = "Provincia de Linares"
$ws.Cells.Item(97,16).Value = 1550
$ws.Cells.Item(97,17).Value = 1
$ws.Cells.Item(97,18).Value = "Hortaliza"

# --- New row 98: Primera ---
$ws.Cells.Item(98,1).Value = 6
$ws.Cells.Item(98,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(98,3).Value = "Metropolitana"
$ws.Cells.Item(98,4).Value = 45225
$ws.Cells.Item(98,5).Value = 13
$ws.Cells.Item(98,6).Value = 300000000
$ws.Cells.Item(98,7).Value = "Espárragos"
$ws.Cells.Item(98,8).Value = "Sin especificar"
$ws.Cells.Item(98,9).Value = "Primera"
$ws.Cells.Item(98,10).Value = 2800
$ws.Cells.Item(98,11).Value = 1100
$ws.Cells.Item(98,12).Value = 1300
$ws.Cells.Item(98,13).Value = 1207
$ws.Cells.Item(98,14).Value = "$/kilo"
$ws.Cells.Item(98,15).Value = "Provincia de Linares"
$ws.Cells.Item(98,16).Value = 1207
$ws.Cells.Item(98,17).Value = 1
$ws.Cells.Item(98,18).Value = "Hortaliza"

# --- New row 99: Segunda ---
$ws.Cells.Item(99,1).Value = 6
$ws.Cells.Item(99,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(99,3).Value = "Metropolitana"
$ws.Cells.Item(99,4).Value = 45225
$ws.Cells.Item(99,5).Value = 13
$ws.Cells.Item(99,6).Value = 300000000
$ws.Cells.Item(99,7).Value = "Espárragos"
$ws.Cells.Item(99,8).Value = "Sin especificar"
$ws.Cells.Item(99,9).Value = "Segunda"
$ws.Cells.Item(99,10).Value = 2500
$ws.Cells.Item(99,11).Value = 900
$ws.Cells.Item(99,12).Value = 1000
$ws.Cells.Item(99,13).Value = 940
$ws.Cells.Item(99,14).Value = "$/kilo"
$ws.Cells.Item(99,15).Value = "Provincia de Linares"
$ws.Cells.Item(99,16).Value = 940
$ws.Cells.Item(99,17).Value = 1
$ws.Cells.Item(99,18).Value = "Hortaliza"
